$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 10:46"

# Row 13 - A Coruña
$ws.Range("B13").Value = 537
$ws.Range("D13").Value = 524
$ws.Range("E13").Value = 13

# Row 17 - Pontevedra
$ws.Range("B17").Value = 448
$ws.Range("D17").Value = 445

# Row 48 - Lugo
$ws.Range("D48").Value = 59
$ws.Range("E48").Value = 3
